$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of data for the new mail log entry ---
$ws = $wb.Worksheets.Item("Logs")

$row = 21
$ws.Cells.Item($row, 1).Value  = "Nieuwe bestelling"
$ws.Cells.Item($row, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 4).Value  = "Inkoop / Bestellingen"
$ws.Cells.Item($row, 6).Value  = "2025-08-28 21:06:14"
$ws.Cells.Item($row, 7).Value  = "Ja"
$ws.Cells.Item($row, 8).Value  = "Nee"
$ws.Cells.Item($row, 9).Value  = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range("$($col)2:$($col)20")
    $newRange = $ws.Range("$($col)2:$($col)21")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Sheet "Dashboard": update the count for "Inkoop / Bestellingen" ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 2
